$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Body paragraph 2: "MaxAir provides an 'Away' option, ..."
#    Split "MaxAir " into a proofErr-wrapped "MaxAir" run + a separate
#    space run.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$xml2 = '<w:p ' + $W + ' xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="17B24648" w14:textId="3286E11F" w:rsidR="00C1132A" w:rsidRDefault="0027538D" w:rsidP="00C1132A">' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>MaxAir</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r w:rsidR="00C1132A"><w:t>provides an ‘Away’ option, which can operate in two different modes, either using a time schedule or not.</w:t></w:r>' + `
    '</w:p>'
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 2) Body paragraph 20 (Note 1): "... is un-ticked." -> "... is unticked."
# ---------------------------------------------------------------------------
$p20 = $d.Paragraphs(20)
$xml20 = '<w:p ' + $W + ' xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6C11FDB2" w14:textId="1A4030F0" w:rsidR="005C5AEC" w:rsidRDefault="005C5AEC" w:rsidP="00C1132A">' + `
    '<w:r><w:t>Note</w:t></w:r>' + `
    '<w:r w:rsidR="00157F38"><w:t xml:space="preserve"> 1</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">: </w:t></w:r>' + `
    '<w:r w:rsidR="00157F38"><w:t>O</w:t></w:r>' + `
    '<w:r><w:t>nly one ‘Away’ schedule can be active at any point in time. Once an ‘Away’ schedule has been activated, the ‘Enable Away Schedule’ checkbox will be disabled</w:t></w:r>' + `
    '<w:r w:rsidR="00782896"><w:t xml:space="preserve"> for any other schedule</w:t></w:r>' + `
    '<w:r w:rsidR="00157F38"><w:t>, until the active ‘Away’ schedule is either deleted or disabled</w:t></w:r>' + `
    '<w:r w:rsidR="00782896"><w:t xml:space="preserve"> or </w:t></w:r>' + `
    '<w:r w:rsidR="00782896"><w:t>‘Enable Away Schedule’ checkbox</w:t></w:r>' + `
    '<w:r w:rsidR="00782896"><w:t xml:space="preserve"> is unticked.</w:t></w:r>' + `
    '</w:p>'
$p20.Range.InsertXML($xml20)

# ---------------------------------------------------------------------------
# 3) Body paragraph 21 (Note 2): split the opening run so "Away " stands on
#    its own as a separate run.
# ---------------------------------------------------------------------------
$p21 = $d.Paragraphs(21)
$xml21 = '<w:p ' + $W + ' xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="16F34F7D" w14:textId="1D31B291" w:rsidR="00157F38" w:rsidRDefault="00157F38" w:rsidP="00C1132A">' + `
    '<w:r><w:t xml:space="preserve">Note 2: When operating in </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Away </w:t></w:r>' + `
    '<w:r><w:t>Scheduled Mode, the system will swap to ‘</w:t></w:r>' + `
    '<w:r w:rsidR="00782896"><w:t>Timer’ Mode irrespective of the Mode set on the ‘Home’ screen</w:t></w:r>' + `
    '<w:r w:rsidR="00F72E62"><w:t xml:space="preserve">, the original mode will be restored when ‘Away’ is deactivated. </w:t></w:r>' + `
    '</w:p>'
$p21.Range.InsertXML($xml21)

# ---------------------------------------------------------------------------
# 4) First-page header: split "MaxAir " the same way as (1), keeping the
#    bookmark wrapping just the "MaxAir" + space runs.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$hp = $hdr.Range.Paragraphs.Item(1)
$xmlH = '<w:p ' + $W + ' xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0CA57FAB" w14:textId="6E037F46" w:rsidR="00C8453E" w:rsidRDefault="00C8453E" w:rsidP="00C8453E">' + `
    '<w:pPr><w:pStyle w:val="Heading1"/><w:jc w:val="center"/></w:pPr>' + `
    '<w:bookmarkStart w:id="0" w:name="_Hlk64485398"/>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>MaxAir</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r w:rsidR="00C1132A"><w:t>– Away Function</w:t></w:r>' + `
    '</w:p>'
$hp.Range.InsertXML($xmlH)

Write-Output "edits applied"
